# Regenerate the "K" column (column G) values for the save_data sheet.
# This mirrors the author's regen of save_data to use K (strikeouts)
# instead of the previous Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New K values keyed by row number (row 2 = first data row ... row 40 = last).
$kValues = @{
    2  = 7
    3  = 6
    4  = 1
    5  = 2
    6  = 6
    7  = 8
    8  = 5
    9  = 6
    10 = 3
    11 = 6
    12 = 2
    13 = 3
    14 = 6
    15 = 11
    16 = 10
    17 = 8
    18 = 1
    19 = 8
    20 = 5
    21 = 6
    22 = 8
    23 = 7
    24 = 10
    25 = 5
    26 = 7
    27 = 4
    28 = 9
    29 = 6
    30 = 9
    31 = 2
    32 = 5
    33 = 11
    34 = 6
    35 = 5
    36 = 7
    37 = 1
    38 = 6
    39 = 6
    40 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
